$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Text content (new / updated) ----
$OBJ_NEW_1 = 'Ensinar a linguagem gráfica normalizada internacionalmente para representação de máquinas e equipamentos que integram os processos de engenharia.Desenvolver o raciocínio espacial e a criatividade de representação.'
$PROG_RESUMIDO_NEW = '-Introdução-Teoria Elementar do Desenho Projetivo-Projeções Ortogonais pelo 1º Diedro-Projeções Ortogonais pelo 3º Diedro-Leitura e Interpretação de Desenhos-Escalas-Desenhos com Instrumentos-Cortes e Representações Convencionais-Projeções Auxiliares-Cotação-Desenhos de Conjuntos e Detalhes-Aplicação de Tolerâncias e Ajustes -Símbolos de Acabamento Superficial-Desenho de Elementos de Máquina-Desenho de Equipamentos e Acessórios'
$PROGRAMA_NEW = '1 - INTRODUÇÃOApresentação e definição da disciplina, destacando a importância do desenho na engenharia; Normas ABNT e ISO.2 - TEORIA ELEMENTAR DO DESENHO PROJETIVORepresentação de vistas como sistema internacional; representação de arestas visíveis e invisíveis; linhas de centro e eixos de simetria.3 - PROJEÇÕES ORTOGONAIS PELO 1º DIEDROPrincípio fundamental; projeções principais; rebatimentos convencionados.4 - PROJEÇÕES ORTOGONAIS PELO 3º DIEDROPrincípio fundamental; projeções principais; rebatimentos convencionados.5 - LEITURA E INTERPRETAÇÃO DE DESENHOSLeitura por meio de esboço em perspectiva e mediante construção de modelos.6 - ESCALASDefinição e normalização7 - DESENHOS COM INSTRUMENTOSRegras para emprego dos esquadros, compasso e régua "T"; disposição do desenho nas folhas padronizadas.8 - CORTES E REPRESENTAÇÕES CONVENCIONAISPrincípios fundamentais; aplicações; tipos normalizados; representações e regras para traçado; seções e rupturas.9 - PROJEÇÕES AUXILIARESPrincípios fundamentais; finalidades e aplicações; representações normalizadas.10 - COTAÇÃORegras de colocação e distribuição de cotas.11 - DESENHOS DE CONJUNTOS E DETALHESDefinições; tipos recomendados de legenda e lista de peça; formas de numeração de desenhos; regras práticas para execução e verificação de desenhos.12 - APLICAÇÃO DE TOLERÂNCIAS E AJUSTESDefinição e finalidades; sistema ISO; uso de tabelas e indicação nos desenhos.13 - SÍMBOLOS DE ACABAMENTO SUPERFICIALDefinição; simbologia normalizada; aplicações.14 - DESENHO DE ELEMENTOS DE MÁQUINADefinições, aplicações, tipos, proporções e representações convencionais de: roscas, parafusos, porcas, arruelas, polias, correias e chavetas.15 - DESENHO DE EQUIPAMENTOS E ACESSÓRIOSDesenho de conjunto e detalhes envolvendo elementos de ligação e de máquinas com aplicação de tabelas e catálogos.'
$AVAL_METODO_NEW = 'A avaliação é continuada e constará de duas provas objetivas (Pi) realizadas ao longo do curso (antes da recuperação), bem como de exercícios práticos realizados em sala de aula e extra classe (TC/TS).'
$CRITERIO_NEW = 'NOTA FINAL = [(MédiaTC/TS)x0,2] + [(MédiaPi)x0,8]'
$NORMA_REC_NEW = '- A recuperação deverá consistir de uma prova englobando a matéria toda do semestre.- A média final (pós-recuperação) deverá ser composta por uma média simples entre a nota do semestre (nota final) e a da prova de recuperação.'
$BIBLIO_NEW = '1 - ABNT - COLETÂNEA DE NORMAS DE DESENHO TÉCNICONormas Técnicas publicadas pela ABNT2 - DESENHO BÁSICO NA ENGENHARIARibeiro, Antonio Clélio - Apostila publicada pela FAENQUIL3 - FUNDAMENTOS DE DIBUJO EM INGENIERIALuzader, Warren J. - Ed. Comp. Editorial Continental - México4 - MANUAL DE DESENHO TÉCNICOManfé, G./ Scarato, G./ Pozza, R. - Ed. Renovada Livros Culturais Ltda.5 - EXPRESSÃO GRÁFICA - DESENHO TÉCNICOHoelsher, R. P./ Springer, C.H./ Dobrovolny, J.S. - Ed. LTC Editora S.A.6 - DESENHO TÉCNICOFrench, Thomas E. - Editora Globo7 - DESENHO TÉCNICOBachmann, A./ Forberg, R - Editora Globo8 - DESENHISTA DE MÁQUINASEscola PRO-TEC'
$MARTINEZ = '5840820 - Gustavo Aristides Santana Martinez'

# ---- Row 10 (Objetivos): replace placeholder content with the real objectives paragraph ----
$ws.Range("B10").Value = $OBJ_NEW_1
$ws.Range("C10").Value = $OBJ_NEW_1

# ---- Insert a new row above row 13 (shifts old rows 13-21 down to 14-22) ----
$ws.Rows.Item(13).Insert()
# the inserted row copied formatting from row 12 into A13; remove that stray cell entirely
$ws.Range("A13").Clear()

# Copy B/C formatting from row 10 (style 2 / style 3) onto the new row 13 B/C cells
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B13").Value = $MARTINEZ
$ws.Range("C13").Value = $MARTINEZ

# ---- Row 14 (Programa resumido): replace "Semestral" placeholder with the summarized program text ----
$ws.Range("B14").Value = $PROG_RESUMIDO_NEW
$ws.Range("C14").Value = $PROG_RESUMIDO_NEW

# ---- Row 16 (Programa): fill in full program text ----
$ws.Range("B16").Value = $PROGRAMA_NEW
$ws.Range("C16").Value = $PROGRAMA_NEW

# ---- Row 19 (Método): fill in evaluation method text ----
$ws.Range("B19").Value = $AVAL_METODO_NEW
$ws.Range("C19").Value = $AVAL_METODO_NEW

# ---- Row 20 (Critério): fill in grading criteria formula text ----
$ws.Range("B20").Value = $CRITERIO_NEW
$ws.Range("C20").Value = $CRITERIO_NEW

# ---- Row 21 (Norma de recuperação): fill in recovery norm text ----
$ws.Range("B21").Value = $NORMA_REC_NEW
$ws.Range("C21").Value = $NORMA_REC_NEW

# ---- Row 22 (Bibliografia): fill in bibliography text ----
$ws.Range("B22").Value = $BIBLIO_NEW
$ws.Range("C22").Value = $BIBLIO_NEW

# ---- Column layout fix: col 1 width rule should only apply to column A now ----
$ws.Columns.Item(1).ColumnWidth = 30.7109375
